$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data values (Fecha, Volumen, Precio minimo/maximo/promedio, Precio $/Kg)
# between row 2 and row 4 - they represent two different weekly records that were
# reordered in the source data (the rest of the row values are identical).

$cols = @("D", "J", "K", "L", "M", "P")

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell4 = $ws.Range($col + "4")
    $tmp = $cell2.Value2
    $cell2.Value2 = $cell4.Value2
    $cell4.Value2 = $tmp
}
